$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.906.70"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "2.239.67"
$ws.Range("E3").Value = "  +2.25%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "270.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.21"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +14.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.636"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.52%  "
$ws.Range("E10").Value = "  +6.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0958"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.34"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +20.03%  "
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.53%  "
$ws.Range("D15").Value = "2.576.70"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("E16").Value = "  +5.59%  "
$ws.Range("D17").Value = "2.249.89"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").Value = "43.898.66"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("E20").Value = "  +4.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("E22").Value = "  -3.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.70"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.68%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.43"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.32%  "
$ws.Range("E27").Value = "  +12.10%  "
$ws.Range("E28").Value = "  +6.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.46"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.80"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0912"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.50"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0353"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.58"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +25.13%  "
$ws.Range("E40").Value = "  +13.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.87"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.57"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.42"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0997"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.42"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("E48").Value = "  +4.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.20"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.445"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("D51").Value = "2.462.21"
$ws.Range("E51").Value = "  +2.38%  "
